$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "CONSTANT COLL EFFIC=0.5"
$ws.Range("B5").Select()
